$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column map: D=4, K=11, L=12, M=13, N=14, O=15, P=16, Q=17, R=18, S=19, T=20

$data = @{
    2  = @{ D=44168; K="Castle Brite"; M=30;  N=16000; O=17000; P=16500; Q="$/caja 16 kilos granel";  S=1031; T=16 }
    3  = @{ D=44552; M=120; N=15500; O=16000; P=15750; Q="$/caja 15 kilos";         S=1050; T=15 }
    4  = @{ D=44189; K="Dina"; L="Primera"; M=80;  P=16562; Q="$/caja 18 kilos";    S=920;  T=18 }
    5  = @{ D=44537; M=60;  N=21000; O=21500; P=21250; Q="$/caja 15 kilos"; R="Región de O'Higgins"; S=1417; T=15 }
    6  = @{ D=44174; M=75;  N=9000;  O=10000; P=9467;  Q="$/caja 10 kilos";         S=947 }
    7  = @{ D=44187; K="Dina"; L="Primera"; M=55;  N=15000; O=16000; P=15455; Q="$/caja 15 kilos granel"; S=1030; T=15 }
    8  = @{ D=44904; M=60;  N=15000; O=16000; P=15500; Q="$/bandeja 10 kilos";      S=1550 }
    9  = @{ D=44904; L="Segunda"; M=30; N=14000; O=14000; P=14000; Q="$/bandeja 10 kilos"; S=1400 }
    10 = @{ D=44165; L="Segunda"; M=60; N=16000; O=17000; P=16500; Q="$/caja 15 kilos granel"; R="Provincia de Limarí"; S=1100 }
    11 = @{ D=44907; L="Primera"; M=120; N=15000; O=16000; P=15500; Q="$/bandeja 10 kilos"; R="Región de O'Higgins"; S=1550; T=10 }
    12 = @{ D=44907; K="Castle Brite"; L="Segunda"; M=60; N=14000; O=14000; P=14000; Q="$/bandeja 10 kilos"; S=1400; T=10 }
    13 = @{ D=44551; M=120; N=15500; O=16000; P=15750; S=1050 }
    14 = @{ D=44176; K="Castle Brite"; M=50; N=17000; O=18000; P=17400; S=967 }
    15 = @{ D=44544; L="Segunda"; M=160; N=16000; O=17000; P=16500; Q="$/caja 15 kilos"; S=1100; T=15 }
    16 = @{ D=44181; K="Modesto"; L="Primera"; M=50; N=20000; O=21000; P=20500; Q="$/caja 18 kilos"; R="Región de Coquimbo"; S=1139; T=18 }
}

$colIndex = @{ D=4; K=11; L=12; M=13; N=14; O=15; P=16; Q=17; R=18; S=19; T=20 }

foreach ($rowNum in $data.Keys) {
    $rowChanges = $data[$rowNum]
    foreach ($col in $rowChanges.Keys) {
        $c = $colIndex[$col]
        $ws.Cells.Item($rowNum, $c).Value = $rowChanges[$col]
    }
}
